# Modificación en plantilla de asignaciones
#
# - Elimina las filas de la "Guía" que describían la columna periodo_pago
#   (periodo / periodo_pago), reacomodando el resto de filas hacia arriba.
# - Elimina la columna "periodo_pago" de las hojas "Horas Extra" y
#   "Asignaciones", reacomodando el resto de columnas hacia la izquierda.
# - Dejar la hoja "Guía" como la hoja activa/seleccionada, con la selección
#   en B16; actualiza la selección recordada en las otras hojas.

$wb = $excel.ActiveWorkbook

$wsGuia  = $wb.Worksheets.Item("Guía")
$wsHoras = $wb.Worksheets.Item("Horas Extra")
$wsAsig  = $wb.Worksheets.Item("Asignaciones")

# Quitar las filas 14 y 15 (periodo / periodo_pago) de la hoja Guía.
[void]$wsGuia.Rows("14:15").Delete()

# Quitar la columna D (periodo_pago) de Horas Extra y Asignaciones.
[void]$wsHoras.Columns("D:D").Delete()
[void]$wsAsig.Columns("D:D").Delete()

# Restaurar/ajustar la selección recordada de cada hoja.
[void]$wsHoras.Activate()
[void]$wsHoras.Range("D2").Select()

[void]$wsAsig.Activate()
[void]$wsAsig.Range("D6").Select()

# La hoja Guía queda como hoja activa, con B16 seleccionado.
[void]$wsGuia.Activate()
[void]$wsGuia.Range("B16").Select()
